# Update the "fidi" watchlist sheet:
#  - refresh the ticker symbols in columns B (Buying Opportunity),
#    C (support Zone), D (long buildup), E (Short buildup) and
#    F (FII ENTERING)
#  - column C grows to 29 entries, pushing the used range down to row 30
#  - the sheet's dimension (A1:F25 -> A1:F30) is recomputed automatically
#    by the engine once row 30 is written

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> 1-based column index
$colIndex = @{ "A" = 1; "B" = 2; "C" = 3; "D" = 4; "E" = 5; "F" = 6 }

# New contents for every data row (row 2 .. row 30), columns B..F.
# $null means "clear this cell" (it had a value before and must become empty).
$rows = @{
    2  = @("NSE:3IINFOLTD",  "NSE:AUTOIND",    "NSE:ALKEM", "NSE:ADANIENSOL", "NSE:COFORGE")
    3  = @("NSE:BALKRISHNA", "NSE:BAJAJHCARE", $null,       "NSE:ADANIPORTS", $null)
    4  = @("NSE:DENORA",     "NSE:BCONCEPTS",  $null,       "NSE:AMBER",      $null)
    5  = @("NSE:GULFPETRO",  "NSE:BHARATGEAR", $null,       "NSE:BDL",        $null)
    6  = @("NSE:HUBTOWN",    "NSE:BLISSGVS",   $null,       "NSE:BHARATFORG", $null)
    7  = @("NSE:IZMO",       "NSE:CYBERMEDIA", $null,       "NSE:BOSCHLTD",   $null)
    8  = @("NSE:JTLIND",     "NSE:DBOL",       $null,       "NSE:HUDCO",      $null)
    9  = @("NSE:LOTUSEYE",   "NSE:DIACABS",    $null,       "NSE:ICICIGI",    $null)
    10 = @("NSE:NAGREEKEXP", "NSE:ENIL",       $null,       "NSE:IRCTC",      $null)
    11 = @("NSE:OSWALAGRO",  "NSE:FAZE3Q",     $null,       "NSE:IRFC",       $null)
    12 = @("NSE:RML",        "NSE:GANDHITUBE", $null,       "NSE:LICI",       $null)
    13 = @("NSE:SALSTEEL",   "NSE:HEIDELBERG", $null,       "NSE:NESTLEIND",  $null)
    14 = @($null,            "NSE:HINDWAREAP", $null,       "NSE:PFC",        $null)
    15 = @($null,            "NSE:INDOTECH",   $null,       $null,            $null)
    16 = @($null,            "NSE:INTLCONV",   $null,       $null,            $null)
    17 = @($null,            "NSE:KIMS",       $null,       $null,            $null)
    18 = @($null,            "NSE:KRSNAA",     $null,       $null,            $null)
    19 = @($null,            "NSE:LAXMICOT",   $null,       $null,            $null)
    20 = @($null,            "NSE:MARKSANS",   $null,       $null,            $null)
    21 = @($null,            "NSE:MON100",     $null,       $null,            $null)
    22 = @($null,            "NSE:NYKAA",      $null,       $null,            $null)
    23 = @($null,            "NSE:ORIENTCER",  $null,       $null,            $null)
    24 = @($null,            "NSE:PDSL",       $null,       $null,            $null)
    25 = @($null,            "NSE:PETRONET",   $null,       $null,            $null)
    26 = @($null,            "NSE:PFC",        $null,       $null,            $null)
    27 = @($null,            "NSE:PITTIENG",   $null,       $null,            $null)
    28 = @($null,            "NSE:PNBGILTS",   $null,       $null,            $null)
    29 = @($null,            "NSE:PNBHOUSING", $null,       $null,            $null)
    30 = @($null,            "NSE:RAYMOND",    $null,       $null,            $null)
}

# Rows 26-30 are brand new -- column A needs the running index (24..28)
# with the same style as the existing index cells (A2:A25), so copy the
# formatting from A25 down before writing the new numbers.
$ws.Range("A25").Copy() | Out-Null
$ws.Range("A26:A30").PasteSpecial(-4122) | Out-Null

$newIndex = @{ 26 = 24; 27 = 25; 28 = 26; 29 = 27; 30 = 28 }
foreach ($r in $newIndex.Keys) {
    $ws.Cells.Item($r, $colIndex["A"]).Value = $newIndex[$r]
}

$columns = @("B", "C", "D", "E", "F")
foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $v = $vals[$i]
        if ($null -eq $v) {
            $ws.Cells.Item($r, $colIndex[$col]).Value = ""
        } else {
            $ws.Cells.Item($r, $colIndex[$col]).Value = $v
        }
    }
}

Write-Output "fidi sheet refreshed"
